$d = $word.ActiveDocument

# 1. Replace "Trung Tâm Cung Ứng Vật Tư – Viễn Thông TP. Hồ Chí Minh" -> "CUVT-HCM"
$d.Content.Find.Execute("Trung Tâm Cung Ứng Vật Tư", $true, $false, $false, $false, $false,
                         $true, 1, $false, "CUVT-HCM", 2)

# Need to remove trailing remnants " – Viễn Thông TP. Hồ Chí Minh" if the above only replaced part
$d.Content.Find.Execute(" – Viễn Thông TP. Hồ Chí Minh", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)

# 2. Replace monetary amount
$d.Content.Find.Execute("19.759.974.500", $true, $false, $false, $false, $false,
                         $true, 1, $false, "19.716.877.500", 2)

# 3. Replace number in words
$d.Content.Find.Execute("mười chín  tỉ bảy trăm năm mươi chín  triệu chín trăm bảy mươi bốn  nghìn năm trăm  ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "mười chín  tỉ bảy trăm mười sáu  triệu tám trăm bảy mươi bảy  nghìn năm trăm  ", 2)
